$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update lose_coef (column O) from 0.95 to 0.93 for rows 10, 11, 12, 14
$ws.Range("O10").Value = 0.93
$ws.Range("O11").Value = 0.93
$ws.Range("O12").Value = 0.93
$ws.Range("O14").Value = 0.93
